$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B201").Value = 7600
$ws.Range("B201").Font.Bold = $true
$ws.Range("C201").Value = 2879095
$ws.Range("D201").Value = 3710409
$ws.Range("E201").Value = 587769
$ws.Range("F201").Value = 117140

$ws.Range("A202").Value = 42625
$ws.Range("B202").Value = 10560
$ws.Range("C202").Value = 3857171
$ws.Range("D202").Value = 4504189
$ws.Range("E202").Value = 609200
$ws.Range("F202").Value = 341200

$ws.Range("A203").Value = 42540
$ws.Range("B203").Value = 10200
$ws.Range("C203").Value = 4074000
$ws.Range("D203").Value = 5415300
$ws.Range("E203").Value = 351000
$ws.Range("F203").Value = 828500

$ws.Range("A204").Value = 42604
$ws.Range("B204").Value = 10200
$ws.Range("B204").Font.Bold = $true
$ws.Range("C204").Value = 3881171
$ws.Range("D204").Value = 4521724
$ws.Range("E204").Value = 628200
$ws.Range("F204").Value = 354020

$ws.Range("A205").Value = 42604
$ws.Range("B205").Value = 7200
$ws.Range("B205").Font.Bold = $true
$ws.Range("C205").Value = 1275200
$ws.Range("D205").Value = 1872700
$ws.Range("E205").Value = 246100
$ws.Range("F205").Value = 236700

$ws.Range("A206").Value = 42598
$ws.Range("B206").Value = 10800
$ws.Range("C206").Value = 2845280
$ws.Range("D206").Value = 5236719
$ws.Range("E206").Value = 560386
$ws.Range("F206").Value = 107690

$ws.Range("A207").Value = 42598
$ws.Range("B207").Value = 8900
$ws.Range("C207").Value = 3881171
$ws.Range("D207").Value = 4521724
$ws.Range("E207").Value = 628200
$ws.Range("F207").Value = 982220

$ws.Range("A208").Value = 42598
$ws.Range("B208").Value = 8900
$ws.Range("C208").Value = 4149000
$ws.Range("D208").Value = 4866960
$ws.Range("E208").Value = 365000
$ws.Range("F208").Value = 792000

$ws.Range("A209").Value = 42606
$ws.Range("B209").Value = 8900
$ws.Range("B209").Font.Bold = $true
$ws.Range("C209").Formula = "=1030937.47+151750"
$ws.Range("D209").Value = 752300
$ws.Range("E209").Value = 289400
$ws.Range("F209").Formula = "=652940+68000"

$ws.Range("A210").Value = 42604
$ws.Range("B210").Value = 7000
$ws.Range("B210").Font.Bold = $true
$ws.Range("C210").Formula = "=142200+1121560"
$ws.Range("D210").Formula = "=468250+62490"
$ws.Range("E210").Value = 287800
$ws.Range("F210").Formula = "=53000+696940"

$ws.Range("A211").Value = 42587
$ws.Range("B211").Value = 10000
$ws.Range("C211").Value = 3217500
$ws.Range("D211").Value = 4713612
$ws.Range("E211").Value = 854500
$ws.Range("F211").Value = 336695

$ws.Range("A212").Value = 42611
$ws.Range("B212").Value = 10400
$ws.Range("B212").Font.Bold = $true
$ws.Range("C212").Value = 3761405
$ws.Range("D212").Value = 3540378
$ws.Range("E212").Value = 523963
$ws.Range("F212").Value = 107690

$ws.Range("A213").Value = 42559
$ws.Range("B213").Value = 9600
$ws.Range("C213").Value = 2426500
$ws.Range("D213").Value = 4803272
$ws.Range("E213").Value = 931500
$ws.Range("F213").Value = 775421.75

$ws.Range("A214").Value = 42590
$ws.Range("B214").Value = 8900
$ws.Range("B214").Font.Bold = $true
$ws.Range("C214").Value = 1683500
$ws.Range("D214").Formula = "=2260000+268000+25000"
$ws.Range("E214").Value = 416500
$ws.Range("F214").Formula = "=12000+142000+46500"

$ws.Range("A215").Value = 42577
$ws.Range("B215").Value = 8200
$ws.Range("B215").Font.Bold = $true
$ws.Range("C215").Value = 1217366
$ws.Range("D215").Formula = "=1412500+268000+25000"
$ws.Range("E215").Value = 292634
$ws.Range("F215").Formula = "=142000+7500+46500"

$ws.Range("A216").Value = 42590
$ws.Range("B216").Value = 10400
$ws.Range("B216").Font.Bold = $true
$ws.Range("C216").Value = 1665250
$ws.Range("D216").Formula = "=2825000+256000+33500"
$ws.Range("E216").Value = 452500
$ws.Range("F216").Formula = "=15000+196000+415000"

$ws.Range("A217").Value = 42597
$ws.Range("B217").Value = 9730
$ws.Range("C217").Value = 3862000
$ws.Range("D217").Value = 5587500
$ws.Range("E217").Value = 456000
$ws.Range("F217").Value = 973200

$ws.Range("A218").Value = 42597
$ws.Range("B218").Value = 9800
$ws.Range("B218").Font.Bold = $true
$ws.Range("C218").Value = 3752928
$ws.Range("D218").Value = 2434049
$ws.Range("E218").Value = 704200
$ws.Range("F218").Value = 314730

$ws.Range("A219").Value = 42592
$ws.Range("B219").Value = 7300
$ws.Range("B219").Font.Bold = $true
$ws.Range("C219").Formula = "=153200+1030938"
$ws.Range("D219").Formula = "=393850+62490"
$ws.Range("E219").Value = 287400
$ws.Range("F219").Formula = "=71500+652940"

$ws.Range("A220").Value = 42591
$ws.Range("B220").Value = 6800
$ws.Range("B220").Font.Bold = $true
$ws.Range("C220").Value = 1956000
$ws.Range("D220").Value = 4784123
$ws.Range("E220").Value = 460000
$ws.Range("F220").Value = 784635

$ws.Range("A221").Value = 42597
$ws.Range("B221").Value = 9730
$ws.Range("C221").Value = 3862000
$ws.Range("D221").Value = 5587500
$ws.Range("E221").Value = 456000
$ws.Range("F221").Value = 973200

$ws.Range("B219").Select()
